$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.351.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.55%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.74"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.17%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0962"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.703.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.40%  "
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.355.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.306.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  +18.34%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("E31").Value = "  +12.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0696"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  +3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.13%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.81%  "
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.448.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.576.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.55%  "
$ws.Range("E51").Value = "  -1.60%  "
